$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 data
$ws.Range("A8").Value = 7

$ws.Range("B8").Value = 2.1346180555555558
$ws.Cells.Item(8, 2).NumberFormat = $ws.Cells.Item(7, 2).NumberFormat

# Set D8 before C8 so shared strings are appended in the same order as the target file
$ws.Range("D8").Value = "Watched Dragon Ball Z, Youtube videos about Minecraft and read Harry Potter."
$ws.Range("C8").Value = "10 Cosas Estúpidas Que Hiciste en Minecraft!!! #1[https://www.youtube.com/watch?v=KvWjAVHz384] (Audiovisual, Spanish, New):36; Dragon Ball Z (Audiovisual, Japanese, Familiar):38;  Harry Potter book 2 (Text-only, English, Familiar):33;"

# Update selection to mimic final saved state (B9 selected)
$ws.Range("B9").Select()
